$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: "A" + " " + "slide" -> "A " + "slide"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "A slide"

# Table cell "a" + " " + "table" -> "a " + "table"
$tbl = $s.Shapes.Item(3).Table
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "a table"

# TextBox: "Plus" + " " + "an" + " " + "image" -> "Plus " + "an " + "image"
$s.Shapes.Item(6).TextFrame.TextRange.Text = "Plus an image"
